$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2265
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = -42
$ws.Range("H2").Value = -38
$ws.Range("I2").Value = -38
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2629
$ws.Range("L2").Value = 1593
$ws.Range("M2").Value = 1036
$ws.Range("N2").Value = 801
$ws.Range("O2").Value = 235
$ws.Range("P2").Value = 284
$ws.Range("Q2").Value = 154
$ws.Range("R2").Value = 22
$ws.Range("S2").Value = -248
$ws.Range("T2").Value = 25
$ws.Range("U2").Value = 129
$ws.Range("V2").Value = 1136
$ws.Range("W2").Value = -0.01
$ws.Range("X2").Value = -1.68
$ws.Range("Y2").Value = -4.67
$ws.Range("Z2").Value = -1.38
$ws.Range("AA2").Value = 153.81
$ws.Range("AB2").Value = 188.09
$ws.Range("AC2").Value = -68
$ws.Range("AD2").Value = -9.82
$ws.Range("AE2").Value = 1460
$ws.Range("AF2").Value = 0.45
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 56857860

# Row 3
$ws.Range("D3").Value = 1969
$ws.Range("E3").Value = -3
$ws.Range("F3").Value = -3
$ws.Range("G3").Value = -46
$ws.Range("H3").Value = -55
$ws.Range("I3").Value = -38
$ws.Range("J3").Value = -16
$ws.Range("K3").Value = 2599
$ws.Range("L3").Value = 1625
$ws.Range("M3").Value = 975
$ws.Range("N3").Value = 755
$ws.Range("O3").Value = 219
$ws.Range("P3").Value = 284
$ws.Range("Q3").Value = 218
$ws.Range("R3").Value = -87
$ws.Range("S3").Value = -102
$ws.Range("T3").Value = 106
$ws.Range("U3").Value = 112
$ws.Range("V3").Value = 1083
$ws.Range("W3").Value = -0.13
$ws.Range("X3").Value = -2.77
$ws.Range("Y3").Value = -4.93
$ws.Range("Z3").Value = -2.09
$ws.Range("AA3").Value = 166.67
$ws.Range("AB3").Value = 173.69
$ws.Range("AC3").Value = -68
$ws.Range("AD3").Value = -14.96
$ws.Range("AE3").Value = 1376
$ws.Range("AF3").Value = 0.73
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 0.99
$ws.Range("AI3").Value = -14.29
$ws.Range("AJ3").Value = 56857860

# Row 4
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("D4").Value = 1957
$ws.Range("E4").Value = -67
$ws.Range("F4").Value = -67
$ws.Range("G4").Value = -22
$ws.Range("H4").Value = -35
$ws.Range("I4").Value = -6
$ws.Range("J4").Value = -29
$ws.Range("K4").Value = 2567
$ws.Range("L4").Value = 1690
$ws.Range("M4").Value = 876
$ws.Range("N4").Value = 688
$ws.Range("O4").Value = 189
$ws.Range("P4").Value = 284
$ws.Range("Q4").Value = -6
$ws.Range("R4").Value = -3
$ws.Range("S4").Value = -6
$ws.Range("T4").Value = 128
$ws.Range("U4").Value = -134
$ws.Range("V4").Value = 1128
$ws.Range("W4").Value = -3.43
$ws.Range("X4").Value = -1.81
$ws.Range("Y4").Value = -0.85
$ws.Range("Z4").Value = -1.37
$ws.Range("AA4").Value = 192.88
$ws.Range("AB4").Value = 168.76
$ws.Range("AC4").Value = -11
$ws.Range("AD4").Value = -84.05
$ws.Range("AE4").Value = 1253
$ws.Range("AF4").Value = 0.73
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 56857860

# Row 5
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("D5").Value = 2031
$ws.Range("E5").Value = -50
$ws.Range("F5").Value = -50
$ws.Range("G5").Value = -93
$ws.Range("H5").Value = -87
$ws.Range("I5").Value = -89
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 2934
$ws.Range("L5").Value = 1856
$ws.Range("M5").Value = 1078
$ws.Range("N5").Value = 856
$ws.Range("O5").Value = 221
$ws.Range("P5").Value = 284
$ws.Range("Q5").Value = -129
$ws.Range("R5").Value = -58
$ws.Range("S5").Value = 83
$ws.Range("T5").Value = 118
$ws.Range("U5").Value = -247
$ws.Range("V5").Value = 1250
$ws.Range("W5").Value = -2.44
$ws.Range("X5").Value = -4.29
$ws.Range("Y5").Value = -11.56
$ws.Range("Z5").Value = -3.17
$ws.Range("AA5").Value = 172.23
$ws.Range("AB5").Value = 137.78
$ws.Range("AC5").Value = -157
$ws.Range("AD5").Value = -4.03
$ws.Range("AE5").Value = 1560
$ws.Range("AF5").Value = 0.41
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 56857860

# Row 6
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()
$ws.Range("D6").Value = 2202
$ws.Range("E6").Value = -107
$ws.Range("F6").Value = -107
$ws.Range("G6").Value = -122
$ws.Range("H6").Value = -130
$ws.Range("I6").Value = -102
$ws.Range("K6").Value = 3065
$ws.Range("L6").Value = 2116
$ws.Range("M6").Value = 948
$ws.Range("N6").Value = 755
$ws.Range("P6").Value = 284
$ws.Range("Q6").Value = -212
$ws.Range("R6").Value = 35
$ws.Range("S6").Value = 216
$ws.Range("T6").Value = 28
$ws.Range("U6").Value = -240
$ws.Range("V6").Value = 1521
$ws.Range("W6").Value = -4.85
$ws.Range("X6").Value = -5.92
$ws.Range("Y6").Value = -12.64
$ws.Range("Z6").Value = -4.34
$ws.Range("AA6").Value = 223.14
$ws.Range("AB6").Value = 102.27
$ws.Range("AC6").Value = -179
$ws.Range("AD6").Value = -3.54
$ws.Range("AE6").Value = 1376
$ws.Range("AF6").Value = 0.46
$ws.Range("AJ6").Value = 56857860

# Rows 7-9: clear all data columns, keep only A/B/C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
